$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.051828742027283
$ws.Range("B1").Value = 2.264790773391724
$ws.Range("C1").Value = 9.462464332580566
$ws.Range("D1").Value = 2.371566772460938
$ws.Range("E1").Value = 1.334019303321838
